# Update countries & provincias Spain
#
# The source COVID-19 dataset was refreshed: most countries simply get
# updated totals, but two pairs of neighbouring countries swapped rank
# order (their row keeps the same rank/position column, but the country
# name + stats that occupy that row change):
#   - rank 55/56: Portugal overtakes Etiopia
#   - rank 147/148/149: Birmania overtakes Guadalupe and Guinea-Bisau
# The "last updated" timestamp in A1 is also bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / last-updated timestamp
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 11 de Septiembre de 2020 a las 16:31"

# Row 5 - India
$ws.Cells.Item(5, 2).Value = 4592533
$ws.Cells.Item(5, 3).Value = 32808
$ws.Cells.Item(5, 4).Value = 3565949
$ws.Cells.Item(5, 5).Value = 950038
$ws.Cells.Item(5, 7).Value = 242
$ws.Cells.Item(5, 8).Value = 76546

# Row 13 - Argentina
$ws.Cells.Item(13, 4).Value = 400121
$ws.Cells.Item(13, 5).Value = 113083
$ws.Cells.Item(13, 7).Value = 87
$ws.Cells.Item(13, 8).Value = 10994

# Row 14 - Chile
$ws.Cells.Item(14, 2).Value = 430535
$ws.Cells.Item(14, 3).Value = 1860
$ws.Cells.Item(14, 4).Value = 403064
$ws.Cells.Item(14, 5).Value = 15621
$ws.Cells.Item(14, 7).Value = 69
$ws.Cells.Item(14, 8).Value = 11850

# Row 24 - Alemania
$ws.Cells.Item(24, 2).Value = 258769
$ws.Cells.Item(24, 3).Value = 662
$ws.Cells.Item(24, 5).Value = 16048
$ws.Cells.Item(24, 7).Value = 2
$ws.Cells.Item(24, 8).Value = 9421

# Row 51 - Portugal moves up past Etiopia, takes the rank-55 row with new stats
$ws.Cells.Item(51, 1).Value = "Portugal"
$ws.Cells.Item(51, 2).Value = 62813
$ws.Cells.Item(51, 3).Value = 687
$ws.Cells.Item(51, 4).Value = 43644
$ws.Cells.Item(51, 5).Value = 17314
$ws.Cells.Item(51, 7).Value = 3
$ws.Cells.Item(51, 8).Value = 1855

# Row 52 - Etiopia drops to the rank-56 row, keeping its previous stats
$ws.Cells.Item(52, 1).Value = "Etiopia"
$ws.Cells.Item(52, 2).Value = 62578
$ws.Cells.Item(52, 4).Value = 23640
$ws.Cells.Item(52, 5).Value = 37964
$ws.Cells.Item(52, 8).Value = 974

# Row 64 - Kirguistan
$ws.Cells.Item(64, 2).Value = 44761
$ws.Cells.Item(64, 3).Value = 77
$ws.Cells.Item(64, 4).Value = 40631
$ws.Cells.Item(64, 5).Value = 3067
$ws.Cells.Item(64, 7).Value = 2
$ws.Cells.Item(64, 8).Value = 1063

# Row 70 - Serbia
$ws.Cells.Item(70, 2).Value = 32228
$ws.Cells.Item(70, 3).Value = 92
$ws.Cells.Item(70, 4).Value = 31100
$ws.Cells.Item(70, 5).Value = 398
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = 730

# Row 77 - Bosnia y Herzegovina
$ws.Cells.Item(77, 2).Value = 22834
$ws.Cells.Item(77, 3).Value = 290
$ws.Cells.Item(77, 4).Value = 15672
$ws.Cells.Item(77, 5).Value = 6476
$ws.Cells.Item(77, 7).Value = 6
$ws.Cells.Item(77, 8).Value = 686

# Row 86 - Republica de Macedonia
$ws.Cells.Item(86, 2).Value = 15555
$ws.Cells.Item(86, 3).Value = 141
$ws.Cells.Item(86, 4).Value = 12994
$ws.Cells.Item(86, 5).Value = 1919
$ws.Cells.Item(86, 7).Value = 5
$ws.Cells.Item(86, 8).Value = 642

# Row 92 - Noruega
$ws.Cells.Item(92, 2).Value = 11924
$ws.Cells.Item(92, 3).Value = 58
$ws.Cells.Item(92, 5).Value = 1288

# Row 143 - Birmania moves up past Guadalupe & Guinea-Bisau, takes the rank-147 row with new stats
$ws.Cells.Item(143, 1).Value = "Birmania"
$ws.Cells.Item(143, 2).Value = 2422
$ws.Cells.Item(143, 3).Value = 272
$ws.Cells.Item(143, 4).Value = 625
$ws.Cells.Item(143, 5).Value = 1783
$ws.Cells.Item(143, 8).Value = 14

# Row 144 - Guadalupe drops to the rank-148 row, keeping its previous stats
$ws.Cells.Item(144, 1).Value = "Guadalupe"
$ws.Cells.Item(144, 2).Value = 2287
$ws.Cells.Item(144, 4).Value = 336
$ws.Cells.Item(144, 5).Value = 1928
$ws.Cells.Item(144, 8).Value = 23

# Row 145 - Guinea-Bisau drops to the rank-149 row, keeping its previous stats
$ws.Cells.Item(145, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(145, 2).Value = 2275
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 1127
$ws.Cells.Item(145, 5).Value = 1109
$ws.Cells.Item(145, 8).Value = 39

# Row 179 - Islas Feroe
$ws.Cells.Item(179, 2).Value = 416
$ws.Cells.Item(179, 3).Value = 1
$ws.Cells.Item(179, 5).Value = 6
